# Delete row 713 ("「奮闘なくして進歩なし」..." post) from the posts sheet.
# This removes the entire row and shifts all following rows up by one,
# which also updates the used range from A1:C815 to A1:C814.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(713).Delete()
